# Add 2022-Q3 data
# ------------------------------------------------------------------
# Target state after edit:
#   sheets (tab order): 总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4
#   - "总计" summary sheet gets a new leading row for 2022-Q3 and all
#     subsequent rows shift down by one.
#   - A brand new "2022-Q3" sheet is inserted (it is a copy of the
#     current "2022-Q2" sheet's layout/data, then the current
#     "2022-Q2" sheet becomes the historical copy retaining the old
#     values, while the former sheet object is renamed to 2022-Q3 and
#     updated with the new fund figures).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Update the "总计" (totals) sheet: insert the 2022-Q3 row at the
#    top of the data and push the existing rows down by one.
# ---------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)

# Make room: clone the formatting of row 4 into the new row 5 first
# so the new row picks up the same cell style as the others.
$totals.Range("A4").Copy($totals.Range("A5"))

# Row 5 <= old row 4 data (2021-Q4, count 1, value 0.23), with the
# running index bumped to 3.
$totals.Range("A5").Value = 3
$totals.Range("B5").Value = "2021-Q4"
$totals.Range("C5").Value = 1
$totals.Range("D5").Value = 0.23

# Row 4 <= old row 3 data (2022-Q1, count 8, value 0.1), index stays 2.
$totals.Range("B4").Value = "2022-Q1"
$totals.Range("C4").Value = 8
$totals.Range("D4").Value = 0.1

# Row 3 <= old row 2 data (2022-Q2, count 1, value 0.35), index stays 1.
$totals.Range("B3").Value = "2022-Q2"
$totals.Range("C3").Value = 1
$totals.Range("D3").Value = 0.35

# Row 2 <= new 2022-Q3 data, index stays 0.
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 1
$totals.Range("D2").Value = 0.38

# ---------------------------------------------------------------
# 2) Duplicate the current "2022-Q2" detail sheet so the historical
#    values are preserved on their own tab, then turn the original
#    sheet into the new "2022-Q3" tab with refreshed figures.
# ---------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item(2)
$q2Sheet.Copy($null, $q2Sheet)

$q2Sheet.Name = "2022-Q3"
$q2Copy = $wb.Worksheets.Item(3)
$q2Copy.Name = "2022-Q2"

# Update the figures on the (now) "2022-Q3" sheet. Columns D-G store
# numeric-looking values as text in this workbook, so force text
# formatting before assigning them to avoid Excel auto-converting the
# strings to numbers.
$q2Sheet.Range("D2:G2").NumberFormat = "@"
$q2Sheet.Range("D2").Value = "11.92"
$q2Sheet.Range("E2").Value = "81.96"
$q2Sheet.Range("F2").Value = "3.18"
$q2Sheet.Range("G2").Value = "0.3791"
$q2Sheet.Range("H2").Value = 8

# ---------------------------------------------------------------
# 3) Restore the original selection state: the last tab ("2021-Q4")
#    was the active sheet/cell before the edit.
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Select() | Out-Null
$lastSheet.Range("A1").Select() | Out-Null
